# Rename columns in data files
# Header row currently contains: AirlineId, Name, AirlineCode
# Rename to snake_case: airline_id, name, airline_code

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "airline_id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "airline_code"

$ws.Range("C2").Select()
